# script written for totals ml parameter search
#
# Roster fix for rows 15-17 (Detroit 2023 roster):
#   - James Wiseman's row moves up from row 17 to row 15, and is given his
#     jersey number (13), which was previously blank.
#   - Jared Rhoden (TW) shifts from row 15 down to row 16 (unchanged data).
#   - Buddy Boeheim (TW) shifts from row 16 down to row 17 (unchanged data).
# The "No." index column A (13, 14, 15) and the header row are untouched;
# only columns B..K move with their player.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 15 -> James Wiseman (was row 17's player; gets jersey No. 13)
$ws.Range("B15").Value = 13
$ws.Range("C15").Value = "James Wiseman"
$ws.Range("D15").Value = "C"
$ws.Range("E15").Value = "7-0"
$ws.Range("F15").Value = 240
$ws.Range("G15").Value = "March 31, 2001"
$ws.Range("H15").Value = "us"
$ws.Range("I15").Value = "1"
$ws.Range("J15").Value = "Memphis"
$ws.Range("K15").Value = "https://www.basketball-reference.com/players/w/wisemja01.html"

# Row 16 -> Jared Rhoden (TW) (was row 15's player, data unchanged)
$ws.Range("B16").Value = 8
$ws.Range("C16").Value = "Jared Rhoden (TW)"
$ws.Range("D16").Value = "SG"
$ws.Range("E16").Value = "6-6"
$ws.Range("F16").Value = 210
$ws.Range("G16").Value = "August 27, 1999"
$ws.Range("H16").Value = "us"
$ws.Range("I16").Value = "R"
$ws.Range("J16").Value = "Seton Hall"
$ws.Range("K16").Value = "https://www.basketball-reference.com/players/r/rhodeja01.html"

# Row 17 -> Buddy Boeheim (TW) (was row 16's player, data unchanged)
$ws.Range("B17").Value = 27
$ws.Range("C17").Value = "Buddy Boeheim (TW)"
$ws.Range("D17").Value = "PG"
$ws.Range("E17").Value = "6-6"
$ws.Range("F17").Value = 205
$ws.Range("G17").Value = "November 11, 1999"
$ws.Range("H17").Value = "us"
$ws.Range("I17").Value = "R"
$ws.Range("J17").Value = "Syracuse"
$ws.Range("K17").Value = "https://www.basketball-reference.com/players/b/boehebu01.html"
